$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update current page (bookmark) for "The Spirit of Kaizen" by Robert Maurer
$ws.Range("C14").Value = 39

# Update active selection from B15 to C15
$ws.Range("C15").Select()
